$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("standard")

# Add new column L with header "Tom" and its values
$ws.Range("L1").Value = "Tom"
$ws.Range("L2").Value = 128.84899999999999
$ws.Range("L3").Value = 29.076499999999999
$ws.Range("L4").Value = 27.434799999999999
$ws.Range("L5").Value = 27.4222

# Update selection to match the target state
$ws.Range("F6").Select()
